# Sprint-1 Week-1 : Tests are marked for execution
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the stale "FAIL" result on the existing IssuedForInformation test case (row 5)
$ws.Range("E5").Value = $null

# Add three new test rows for execution (IssuedForReview, RequestForInformation, IssuedForApproval)
# New shared strings are interned in the order the cells are written, so write
# the A-column (FLD_... keys) fully before the B-column (description) values.
$ws.Range("A6").Value = "FLD_Transmittals_New_IssuedForReview"
$ws.Range("A7").Value = "FLD_Transmittals_New_RequestForInformation"
$ws.Range("A8").Value = "FLD_Transmittals_New_IssuedForApproval"

$ws.Range("B6").Value = "Creates a new Transmittal of  Issue Reason  IssuedForReview"
$ws.Range("B7").Value = "Creates a new Transmittal of  Issue Reason  RequestForInformation"
$ws.Range("B8").Value = "Creates a new Transmittal of  Issue Reason  IssuedForApproval"

$ws.Range("C6").Value = "N"
$ws.Range("D6").Value = "Y"
$ws.Range("F6").Value = "Sprint1"

$ws.Range("C7").Value = "N"
$ws.Range("D7").Value = "Y"
$ws.Range("F7").Value = "Sprint1"

$ws.Range("C8").Value = "N"
$ws.Range("D8").Value = "Y"
$ws.Range("F8").Value = "Sprint1"

# Copy styles from row 5 down to the new rows so formatting/borders match
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Extend the data validation ranges to cover the new rows
$ws.Range("C2:D8").Validation.Delete()
$ws.Range("C2:D8").Validation.Add(3, 1, 1, "Y,N")

$ws.Range("F2:F8").Validation.Delete()
$ws.Range("F2:F8").Validation.Add(3, 1, 1, "Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10")

